# fix locator on become partner page
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BecomePartnerPage")

# Update the Marketo checkbox locators to the new form id (139660 -> 144792)
$ws.Range("C19").Value = "#mktoCheckbox_144792_0"
$ws.Range("C20").Value = "#LblmktoCheckbox_144792_0"

# Match the author's final view/selection state in the sheet
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C20").Select()
